$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "# Deployments" column (D) values ---
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 5
$ws.Range("D4").Value = 6
$ws.Range("D5").Value = 3

# --- Add new "Lead Time (Days)" column (E) ---
$ws.Range("E1").Value = "Lead Time (Days)"
$ws.Range("E2").Value = 20
$ws.Range("E3").Value = 24
$ws.Range("E4").Value = 19
$ws.Range("E5").Value = 112

# Match the header formatting used by the other header cells (bold font +
# fill), but with only a left border (this is the last column so it
# shouldn't pick up the interior border of the table header row).
$ws.Range("C1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E1").Borders.LineStyle = -4142  # xlLineStyleNone
$ws.Range("E1").Borders.Item(7).LineStyle = 1  # xlContinuous (left edge)
$ws.Range("E1").Borders.Item(7).Weight = 2  # xlThin

# Widen the new column to fit its header text
$ws.Columns.Item(5).ColumnWidth = 18

# Move the selection to the first data cell of the new column
$ws.Range("E2").Select() | Out-Null
